$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3375
$ws1.Range("F4").Value = 2463
$ws1.Range("F5").Value = 336
$ws1.Range("G5").Value = "已售罄"
$ws1.Range("F6").Value = 344
$ws1.Range("F8").Value = 1103
$ws1.Range("F9").Value = 309
$ws1.Range("F10").Value = 520
$ws1.Range("F14").Value = 549
$ws1.Range("F15").Value = 8646
$ws1.Range("F18").Value = 258
$ws1.Range("F22").Value = 592
$ws1.Range("F24").Value = 1156
$ws1.Range("F26").Value = 2028
$ws1.Range("F27").Value = 2068
$ws1.Range("F29").Value = 1763
$ws1.Range("F33").Value = 34
$ws1.Range("F34").Value = 49
$ws1.Range("F35").Value = 93
$ws1.Range("F36").Value = 186
$ws1.Range("F41").Value = 430
$ws1.Range("F42").Value = 688

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3375
$ws4.Range("F4").Value = 2463
$ws4.Range("F5").Value = 336
$ws4.Range("G5").Value = "已售罄"
$ws4.Range("F6").Value = 344
$ws4.Range("F9").Value = 1103
$ws4.Range("F10").Value = 309
$ws4.Range("F11").Value = 520
$ws4.Range("F14").Value = 549
$ws4.Range("F15").Value = 8646
$ws4.Range("F19").Value = 258
$ws4.Range("F23").Value = 592
$ws4.Range("F25").Value = 1156
$ws4.Range("F27").Value = 2028
$ws4.Range("F28").Value = 2068
$ws4.Range("F29").Value = 1763
$ws4.Range("F33").Value = 34
$ws4.Range("F34").Value = 49
$ws4.Range("F35").Value = 93
$ws4.Range("F36").Value = 186
$ws4.Range("F41").Value = 430
$ws4.Range("F46").Value = 688
